# "Lien Github du projet : ……………" -> "Lien Github du projet : "
# followed by a new run holding the GitHub repository URL, matching the
# commit "MAJ Doc explaination for exam" which fills in the previously
# blank "Lien Github du projet" field.

$d = $word.ActiveDocument

# This exact text (trailing 5-dot ellipsis placeholder) only occurs once in
# the document, in the "Lien Github du projet" paragraph, so the Find below
# cannot collide with the sibling "Nom du projet : " or
# "Lien Drive du projet" placeholder runs.
$oldText = " du projet : ……………"
$newText = " du projet : "

$range = $d.Content
$range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# $range now spans exactly the replaced " du projet : " text (Find with
# wdReplaceOne repositions the range over the replacement). Collapse to its
# end point -- right after the colon + space -- and type the URL there so it
# inherits that run's current formatting (sz/szCs = 21).
$range.Collapse(0)
$insertStart = $range.Start
$range.Text = "https://github.com/flow01-lab/CineBooking-BDD-Eval-Train.git"

# Re-seat a range over just the newly inserted URL text and touch its font
# size (same value, set twice) so Word materialises it into its own run
# with its own rPr instead of silently re-merging into the preceding run.
$urlRange = $d.Range($insertStart, $range.End)
$urlRange.Font.Size = 12
$urlRange.Font.Size = 10.5
